$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6849.8
$ws.Range("I32").Value = 6849.8
$ws.Range("K32").Value = 6849.8
$ws.Range("M32").Value = -6523.8
$ws.Range("H40").Value = 4684.2104
$ws.Range("I40").Value = 3065
$ws.Range("J40").Value = 5116
$ws.Range("K40").Value = 3065
$ws.Range("L40").Value = 5116
$ws.Range("M40").Value = -2890
$ws.Range("N40").Value = -5466
$ws.Range("H55").Value = 228.66667
$ws.Range("J55").Value = 307.33334
$ws.Range("L55").Value = 307.33334
$ws.Range("N55").Value = -735.33334
$ws.Range("H62").Value = 6059.769
$ws.Range("I62").Value = 7475.222
$ws.Range("K62").Value = 7475.222
$ws.Range("M62").Value = -6851.222
$ws.Range("H65").Value = 6059.769
$ws.Range("I65").Value = 7475.222
$ws.Range("K65").Value = 37376.11
$ws.Range("M65").Value = -34256.11
$ws.Range("H70").Value = 2088.2
$ws.Range("I70").Value = 1662.3334
$ws.Range("J70").Value = 2372.111
$ws.Range("K70").Value = 4987.0002
$ws.Range("L70").Value = 7116.333
$ws.Range("M70").Value = -4717.0002
$ws.Range("N70").Value = -7656.333
$ws.Range("H73").Value = 2088.2
$ws.Range("I73").Value = 1662.3334
$ws.Range("J73").Value = 2372.111
$ws.Range("K73").Value = 4987.0002
$ws.Range("L73").Value = 7116.333
$ws.Range("M73").Value = -4051.0002
$ws.Range("N73").Value = -8988.332999999999
$ws.Range("H100").Value = 1489.7858
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082
$ws.Range("H111").Value = 8833
$ws.Range("I111").Value = 6500
$ws.Range("J111").Value = 9999.5
$ws.Range("K111").Value = 19500
$ws.Range("L111").Value = 29998.5
$ws.Range("M111").Value = -16433
$ws.Range("N111").Value = -36132.5
$ws.Range("H112").Value = 3817.0908
$ws.Range("J112").Value = 4251.125
$ws.Range("L112").Value = 12753.375
$ws.Range("N112").Value = -14969.375
$ws.Range("H113").Value = 17211.46
$ws.Range("I113").Value = 4460
$ws.Range("J113").Value = 28141.285
$ws.Range("K113").Value = 4460
$ws.Range("L113").Value = 28141.285
$ws.Range("M113").Value = -1206
$ws.Range("N113").Value = -34649.285
$ws.Range("H116").Value = 4471.4287
$ws.Range("I116").Value = 3300
$ws.Range("K116").Value = 3300
$ws.Range("M116").Value = 142
$ws.Range("H132").Value = 2175.3157
$ws.Range("I132").Value = 2148.8823
$ws.Range("K132").Value = 6446.646900000001
$ws.Range("M132").Value = -3916.646900000001
$ws.Range("H137").Value = 1350.8667
$ws.Range("I137").Value = 1233.0714
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 3699.2142
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -1149.2142
$ws.Range("N137").Value = -14100
$ws.Range("H141").Value = 6238.1924
$ws.Range("I141").Value = 5636.4
$ws.Range("K141").Value = 16909.2
$ws.Range("M141").Value = -11729.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7020.8857
$ws.Range("I32").Value = 6536.485
$ws.Range("K32").Value = 6536.485
$ws.Range("M32").Value = -6249.485
$ws.Range("H63").Value = 5178.1665
$ws.Range("I63").Value = 3626.25
$ws.Range("J63").Value = 8282
$ws.Range("K63").Value = 3626.25
$ws.Range("L63").Value = 8282
$ws.Range("M63").Value = -2940.25
$ws.Range("N63").Value = -9654
$ws.Range("H66").Value = 5178.1665
$ws.Range("I66").Value = 3626.25
$ws.Range("J66").Value = 8282
$ws.Range("K66").Value = 18131.25
$ws.Range("L66").Value = 41410
$ws.Range("M66").Value = -14699.25
$ws.Range("N66").Value = -48274
$ws.Range("H122").Value = 2534.75
$ws.Range("I122").Value = 2534.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7604.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5154.25
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1024.3334
$ws.Range("I94").Value = 508.8
$ws.Range("K94").Value = 508.8
$ws.Range("M94").Value = -57.80000000000001
$ws.Range("H107").Value = 1917.909
$ws.Range("I107").Value = 2101.4285
$ws.Range("K107").Value = 2101.4285
$ws.Range("M107").Value = -181.4285
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040
$ws.Range("H134").Value = 3765.9387
$ws.Range("I134").Value = 3305.152
$ws.Range("K134").Value = 9915.456
$ws.Range("M134").Value = -7380.456

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8834.843999999999
$ws.Range("I31").Value = 11399.2
$ws.Range("J31").Value = 6572.1763
$ws.Range("K31").Value = 11399.2
$ws.Range("L31").Value = 6572.1763
$ws.Range("M31").Value = -11104.2
$ws.Range("N31").Value = -7162.1763
$ws.Range("H34").Value = 8834.843999999999
$ws.Range("I34").Value = 11399.2
$ws.Range("J34").Value = 6572.1763
$ws.Range("K34").Value = 11399.2
$ws.Range("L34").Value = 6572.1763
$ws.Range("M34").Value = -11197.2
$ws.Range("N34").Value = -6976.1763
$ws.Range("H118").Value = 80000
$ws.Range("J118").Value = 80000
$ws.Range("L118").Value = 80000
$ws.Range("N118").Value = -83314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 216.66667
$ws.Range("I92").Value = 225
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 675
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 573
$ws.Range("N92").Value = -3096
$ws.Range("H107").Value = 422.83334
$ws.Range("I107").Value = 479.33334
$ws.Range("J107").Value = 366.33334
$ws.Range("K107").Value = 1438.00002
$ws.Range("L107").Value = 1099.00002
$ws.Range("M107").Value = 481.9999800000001
$ws.Range("N107").Value = -4939.000019999999
$ws.Range("H129").Value = 13779399
$ws.Range("I129").Value = 33433972
$ws.Range("J129").Value = 676350.6
$ws.Range("K129").Value = 100301916
$ws.Range("L129").Value = 2029051.8
$ws.Range("M129").Value = -100296916
$ws.Range("N129").Value = -2039051.8
$ws.Range("H131").Value = 13160252
$ws.Range("I131").Value = 100000920
$ws.Range("J131").Value = 2575.2424
$ws.Range("K131").Value = 300002760
$ws.Range("L131").Value = 7725.7272
$ws.Range("M131").Value = -299997720
$ws.Range("N131").Value = -17805.7272
$ws.Range("H137").Value = 4429.1816
$ws.Range("J137").Value = 10672.333
$ws.Range("L137").Value = 32016.999
$ws.Range("N137").Value = -42216.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4264.8
$ws.Range("I14").Value = 6125
$ws.Range("K14").Value = 6125
$ws.Range("M14").Value = -5957
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 925.0833
$ws.Range("I107").Value = 517.3333
$ws.Range("J107").Value = 1061
$ws.Range("K107").Value = 517.3333
$ws.Range("L107").Value = 1061
$ws.Range("M107").Value = 1402.6667
$ws.Range("N107").Value = -4901
$ws.Range("H113").Value = 502162.12
$ws.Range("I113").Value = 1000799.25
$ws.Range("K113").Value = 1000799.25
$ws.Range("M113").Value = -998629.25
$ws.Range("H126").Value = 3999.8572
$ws.Range("J126").Value = 5499.8335
$ws.Range("L126").Value = 16499.5005
$ws.Range("N126").Value = -21439.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2887.375
$ws.Range("J22").Value = 1800
$ws.Range("L22").Value = 1800
$ws.Range("N22").Value = -2390
$ws.Range("H27").Value = 2887.375
$ws.Range("J27").Value = 1800
$ws.Range("L27").Value = 1800
$ws.Range("N27").Value = -2014
$ws.Range("H40").Value = 4481.3335
$ws.Range("J40").Value = 4857.143
$ws.Range("L40").Value = 4857.143
$ws.Range("N40").Value = -5129.143
$ws.Range("H132").Value = 5705.759
$ws.Range("I132").Value = 6075.875
$ws.Range("K132").Value = 18227.625
$ws.Range("M132").Value = -15697.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 9298.333000000001
$ws.Range("I49").Value = 9497.5
$ws.Range("J49").Value = 8900
$ws.Range("K49").Value = 9497.5
$ws.Range("L49").Value = 8900
$ws.Range("M49").Value = -9267.5
$ws.Range("N49").Value = -9360
$ws.Range("H107").Value = 1530.75
$ws.Range("I107").Value = 842
$ws.Range("J107").Value = 2219.5
$ws.Range("K107").Value = 2526
$ws.Range("L107").Value = 6658.5
$ws.Range("M107").Value = -606
$ws.Range("N107").Value = -10498.5
$ws.Range("H122").Value = 3246.2222
$ws.Range("I122").Value = 2258.2778
$ws.Range("J122").Value = 5222.1113
$ws.Range("K122").Value = 6774.8334
$ws.Range("L122").Value = 15666.3339
$ws.Range("M122").Value = -4324.8334
$ws.Range("N122").Value = -20566.3339
$ws.Range("H126").Value = 10525.5
$ws.Range("J126").Value = 13700.8
$ws.Range("L126").Value = 41102.39999999999
$ws.Range("N126").Value = -46042.39999999999
$ws.Range("H132").Value = 5290.697
$ws.Range("I132").Value = 4612.231
$ws.Range("K132").Value = 13836.693
$ws.Range("M132").Value = -11306.693
